$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("create_cashflow")
$ws2 = $wb.Worksheets.Item("moneymarket_rate")

# ---------------------------------------------------------------------------
# Sheet "create_cashflow": fix the maturity-date formula so it strips the
# trailing unit letter (Y/M/W) with LEN(..)-1 instead of hard-coding a
# single character, and tidy up the IF() argument spacing.
# ---------------------------------------------------------------------------
$ws1.Range("C4").Formula = "=B4+IF(E4=1,360,0)*LEFT(A4,LEN(A4)-1)+IF(F4=1,30,0)*LEFT(A4,LEN(A4)-1)+IF(G4=1,7,0)*LEFT(A4,LEN(A4)-1)"
$ws1.Range("C5:C9").Formula = "=B5+IF(E5=1,360,0)*LEFT(A5,LEN(A5)-1)+IF(F5=1,30,0)*LEFT(A5,LEN(A5)-1)+IF(G5=1,7,0)*LEFT(A5,LEN(A5)-1)"
$ws1.Range("C10").Formula = "=B10+IF(E10=1,360,0)*LEFT(A10,LEN(A10)-1)+IF(F10=1,30,0)*LEFT(A10,LEN(A10)-1)+IF(G10=1,7,0)*LEFT(A10,LEN(A10)-1)"

# ---------------------------------------------------------------------------
# Sheet "moneymarket_rate": refresh the quoted rates (column D) and add a new
# column E implementing get_DF, the money-market discount-factor bootstrap:
#   DF(1) = 1 / (1 + (C-B)/360 * rate)
#   DF(n) = DF(previous) / (1 + (C-B)/360 * rate)   [chained off the T/N DF]
# ---------------------------------------------------------------------------
$ws2.Range("D2").Value2 = 0.014348
$ws2.Range("D3").Value2 = 0.014348
$ws2.Range("D4").Value2 = 0.014876
$ws2.Range("D5").Value2 = 0.015
$ws2.Range("D6").Value2 = 0.01563
$ws2.Range("D7").Value2 = 0.01616
$ws2.Range("D8").Value2 = 0.01685
$ws2.Range("D9").Value2 = 0.01833
$ws2.Range("D10").Value2 = 0.021

# The 12M maturity date follows the corrected formula logic on the other
# sheet (B10 + 30 * LEFT("12M", LEN("12M")-1) = 43094 + 30*12).
$ws2.Range("C10").Value2 = 43454

$ws2.Range("E2").Formula = "=1/(1+(C2-B2)/360*D2)"
$ws2.Range("E3").Formula = "=E2/(1+(C3-B3)/360*D3)"
$ws2.Range("E4").Formula = "=`$E`$3/(1+(C4-B4)/360*D4)"
$ws2.Range("E5:E10").Formula = "=`$E`$3/(1+(C5-B5)/360*D5)"

# Leave the same cell selected/active as in the authored workbook.
$ws2.Activate()
$ws2.Range("E9").Select()
